$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "unnamed: 1_level_1" column header (B2) is replaced by "total", matching B1
$ws.Range("B2").Value = "total"

# The two blank section-header rows ("situação do domicílio" and
# "grandes regiões e unidades da federação") carried no data and are
# removed; Excel shifts the remaining data rows up to close the gaps.
$ws.Rows("5").Delete()
$ws.Rows("7").Delete()
